$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.96399999999999
$ws.Range("B14").Value = 8.915300000000002
$ws.Range("B21").Value = 5.758699999999993
$ws.Range("B23").Value = 5.627899999999999
$ws.Range("B25").Value = 5.956099999999993
